$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 888, pushing the existing
# rows 888-932 down to 891-935.
$ws.Rows.Item(888).EntireRow.Insert()
$ws.Rows.Item(888).EntireRow.Insert()
$ws.Rows.Item(888).EntireRow.Insert()

# Common (constant) values shared by every row of this block.
$mercadoId = 3
$mercado = "Femacal de La Calera"
$region = "Coquimbo"
$codreg = 5
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria = "Plátano"
$variedad = "Sin especificar"
$unidad = "$/caja 20 kilos"
$origen = "Ecuador"
$kgUnidad = 20

function Set-FilaPlatano($fila, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($fila, 1).Value = $mercadoId
    $ws.Cells.Item($fila, 2).Value = $mercado
    $ws.Cells.Item($fila, 3).Value = $region
    $ws.Cells.Item($fila, 4).Value = $fecha
    $ws.Cells.Item($fila, 5).Value = $codreg
    $ws.Cells.Item($fila, 6).Value = $tipo
    $ws.Cells.Item($fila, 7).Value = $productoId
    $ws.Cells.Item($fila, 8).Value = $producto
    $ws.Cells.Item($fila, 9).Value = $categoriaId
    $ws.Cells.Item($fila, 10).Value = $categoria
    $ws.Cells.Item($fila, 11).Value = $variedad
    $ws.Cells.Item($fila, 12).Value = $calidad
    $ws.Cells.Item($fila, 13).Value = $volumen
    $ws.Cells.Item($fila, 14).Value = $precioMin
    $ws.Cells.Item($fila, 15).Value = $precioMax
    $ws.Cells.Item($fila, 16).Value = $precioProm
    $ws.Cells.Item($fila, 17).Value = $unidad
    $ws.Cells.Item($fila, 18).Value = $origen
    $ws.Cells.Item($fila, 19).Value = $precioKg
    $ws.Cells.Item($fila, 20).Value = $kgUnidad
}

# New rows of data (same date, three quality grades).
Set-FilaPlatano 888 44753 "Maduro"          200 21000 21000 21000 1050
Set-FilaPlatano 889 44753 "Pintón"          280 22000 22000 22000 1100
Set-FilaPlatano 890 44753 "Primera Pintón"  280 23000 23000 23000 1150
